$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: fill A2:J2 with "Unassigned"
$ws.Range("A2:J2").Value = "Unassigned"

# Update row 3: only A3 has "Microstomus kitt"
$ws.Range("A3").Value = "Microstomus kitt"

# Delete the old rows 4, 5, 6 content (shift cells up by deleting entire rows)
$ws.Range("A4:J6").EntireRow.Delete()
